$wb = $excel.ActiveWorkbook

# Rename test-data person "Julie Carthane" -> "Aja Mount" everywhere it
# appears (Users sheet roster + the matching sample row on the GiftLog
# sheet both point at the same shared-string value in the workbook).
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Aja Mount"

$wsGiftLog = $wb.Worksheets.Item("GiftLog")
$wsGiftLog.Range("B2").Value = "Aja Mount"

# Leave the pre-existing selection on GiftLog parked at B2 (previously the
# active tab, was sitting on C18).
$wsGiftLog.Range("B2").Select()

# Users becomes the active tab now, with the cursor on F10 (previously
# AppName/ModuleName/GiftLog.. ended up with Users active at A2).
$wsUsers.Activate()
$wsUsers.Range("F10").Select()
